# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Leve profit tables
# across all class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 56.642857
$ws.Range("I11").Value = 56.642857
$ws.Range("K11").Value = 56.642857
$ws.Range("M11").Value = 83.35714300000001
$ws.Range("H40").Value = 6274.727
$ws.Range("J40").Value = 8361.75
$ws.Range("L40").Value = 8361.75
$ws.Range("N40").Value = -8711.75
$ws.Range("H116").Value = 12926.637
$ws.Range("I116").Value = 13424.125
$ws.Range("K116").Value = 13424.125
$ws.Range("M116").Value = -9982.125
$ws.Range("H132").Value = 3348.4055
$ws.Range("I132").Value = 1743.8667
$ws.Range("J132").Value = 10225
$ws.Range("K132").Value = 5231.6001
$ws.Range("L132").Value = 30675
$ws.Range("M132").Value = -2701.6001
$ws.Range("N132").Value = -35735

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2859443.5
$ws.Range("I32").Value = 268.0357
$ws.Range("K32").Value = 268.0357
$ws.Range("M32").Value = 18.96429999999998
$ws.Range("H61").Value = 3811.9473
$ws.Range("J61").Value = 6000
$ws.Range("L61").Value = 6000
$ws.Range("N61").Value = -6424
$ws.Range("H97").Value = 1185.1177
$ws.Range("I97").Value = 1137.9166
$ws.Range("J97").Value = 1298.4
$ws.Range("K97").Value = 1137.9166
$ws.Range("L97").Value = 1298.4
$ws.Range("M97").Value = -641.9166
$ws.Range("N97").Value = -2290.4
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()  # was -49843
$ws.Range("H136").Value = 3811.9473
$ws.Range("J136").Value = 6000
$ws.Range("L136").Value = 18000
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()  # was 73
$ws.Range("H33").Value = 9975
$ws.Range("I33").Value = 9975
$ws.Range("K33").Value = 9975
$ws.Range("M33").Value = -9639
$ws.Range("H99").Value = 111112020
$ws.Range("I99").Value = 125000930
$ws.Range("K99").Value = 125000930
$ws.Range("M99").Value = -124999432
$ws.Range("H134").Value = 3390.7144
$ws.Range("I134").Value = 939.25
$ws.Range("J134").Value = 18099.5
$ws.Range("K134").Value = 2817.75
$ws.Range("L134").Value = 54298.5
$ws.Range("M134").Value = -282.75
$ws.Range("N134").Value = -59368.5
$ws.Range("H140").Value = 114472.25
$ws.Range("J140").Value = 114472.25
$ws.Range("L140").Value = 114472.25
$ws.Range("N140").Value = -124832.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 1304.1
$ws.Range("I32").Value = 1005.25
$ws.Range("K32").Value = 1005.25
$ws.Range("M32").Value = -689.25
$ws.Range("H33").Value = 3645.182
$ws.Range("I33").Value = 1774.625
$ws.Range("K33").Value = 1774.625
$ws.Range("M33").Value = -1395.625
$ws.Range("H38").Value = 11001.875
$ws.Range("I38").Value = 3002.5
$ws.Range("J38").Value = 35000
$ws.Range("K38").Value = 3002.5
$ws.Range("L38").Value = 35000
$ws.Range("M38").Value = -2625.5
$ws.Range("N38").Value = -35754
$ws.Range("H46").Value = 11001.875
$ws.Range("I46").Value = 3002.5
$ws.Range("J46").Value = 35000
$ws.Range("K46").Value = 3002.5
$ws.Range("L46").Value = 35000
$ws.Range("M46").Value = -2791.5
$ws.Range("N46").Value = -35422
$ws.Range("H86").Value = 5798
$ws.Range("J86").Value = 5697.5
$ws.Range("L86").Value = 5697.5
$ws.Range("N86").Value = -7943.5
$ws.Range("H89").Value = 5798
$ws.Range("J89").Value = 5697.5
$ws.Range("L89").Value = 28487.5
$ws.Range("N89").Value = -39719.5
$ws.Range("H105").Value = 1700
$ws.Range("I105").Value = 1700
$ws.Range("K105").Value = 1700
$ws.Range("M105").Value = 47
$ws.Range("H122").Value = 650.8
$ws.Range("I122").Value = 650.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1952.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 497.6000000000001
$ws.Range("N122").ClearContents()  # was -6089.000019999999
$ws.Range("H132").Value = 2687.2727
$ws.Range("I132").Value = 2453.9473
$ws.Range("K132").Value = 7361.841899999999
$ws.Range("M132").Value = -4831.841899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 2
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()  # was -3000576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -2058
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()  # was -51372
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()  # was -156864
$ws.Range("H113").Value = 8488.111000000001
$ws.Range("I113").Value = 5464.3335
$ws.Range("K113").Value = 5464.3335
$ws.Range("M113").Value = -3294.3335
$ws.Range("H134").Value = 106108
$ws.Range("J134").Value = 106108
$ws.Range("L134").Value = 318324
$ws.Range("N134").Value = -323394
$ws.Range("H135").Value = 233333.33
$ws.Range("J135").Value = 233333.33
$ws.Range("L135").Value = 233333.33
$ws.Range("N135").Value = -243473.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7499.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 7499.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 7499.5
$ws.Range("M7").ClearContents()  # was -2598.75
$ws.Range("N7").Value = -7723.5
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820
$ws.Range("H126").Value = 7499.5
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 7499.5
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 22498.5
$ws.Range("M126").ClearContents()  # was -5662.25
$ws.Range("N126").Value = -27438.5
$ws.Range("H136").Value = 4333.3335
$ws.Range("I136").Value = 4333.3335
$ws.Range("K136").Value = 13000.0005
$ws.Range("M136").Value = -10450.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4532.4707
$ws.Range("I126").Value = 2130.2
$ws.Range("K126").Value = 6390.599999999999
$ws.Range("M126").Value = -3920.599999999999
$ws.Range("H136").Value = 2992.0625
$ws.Range("I136").Value = 2232.375
$ws.Range("J136").Value = 5271.125
$ws.Range("K136").Value = 6697.125
$ws.Range("L136").Value = 15813.375
$ws.Range("M136").Value = -4147.125
$ws.Range("N136").Value = -20913.375

Write-Output "Applied scheduled Leve profit updates across all sheets."
